$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.13
$ws.Range("E2").Value = 0.235
$ws.Range("G2").Value = 0.05996287128712871
$ws.Range("H2").Value = 0.05996287128712871
$ws.Range("I2").Value = 0.07608258823498562
$ws.Range("J2").Value = 0.05935115370101465
$ws.Range("K2").Value = 28.58
$ws.Range("L2").Value = 0.05895214521452145
$ws.Range("M2").Value = 12.0772
$ws.Range("N2").Value = 0.04380558578164672
$ws.Range("O2").Value = 0.4225752274317705
$ws.Range("P2").Value = 12.0772
$ws.Range("Q2").Value = 0.04380558578164672
$ws.Range("R2").Value = 0.4225752274317705
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 17.44
$ws.Range("V2").Value = 0.06325716358360536
$ws.Range("W2").Value = 0.1230889235569423
$ws.Range("X2").Value = 0.06345807899599111
$ws.Range("Y2").Value = 0.05963084456095118
$ws.Range("Z2").Value = 2.086886394824474
$ws.Range("AA2").Value = 0.1210294192434357
$ws.Range("AB2").Value = 0.06331852767618502
$ws.Range("AC2").Value = 0.05771089156725072
$ws.Range("AD2").Value = 10.3
$ws.Range("AE2").Value = 1.175806118394853
$ws.Range("AF2").Value = 11.47580611839485
$ws.Range("AG2").Value = -5.964193881605144
$ws.Range("AH2").Value = 0.03996090852327472
$ws.Range("AI2").Value = 0.04373804996797797
$ws.Range("AJ2").Value = -0.02211124272832835
$ws.Range("AK2").Value = -0.02435002858962248
$ws.Range("AL2").Value = 0.307
$ws.Range("AM2").Value = 0.114
$ws.Range("AN2").Value = 0.2529469548133595
$ws.Range("AO2").Value = 120.5211726384365
$ws.Range("AP2").Value = -0.1464684155600477
$ws.Range("AQ2").Value = 324.5614035087719

# Row 3
$ws.Range("D3").Value = 0.0824
$ws.Range("E3").Value = 0.132
$ws.Range("G3").Value = 0.08353221957040573
$ws.Range("H3").Value = 0.08353221957040573
$ws.Range("I3").Value = 0.0803500397772474
$ws.Range("J3").Value = 0.06410129563624692
$ws.Range("K3").Value = 7.89
$ws.Range("L3").Value = 0.06276849642004773
$ws.Range("M3").Value = 4.92
$ws.Range("N3").Value = 0.05956416464891041
$ws.Range("O3").Value = 0.623574144486692
$ws.Range("P3").Value = 4.92
$ws.Range("Q3").Value = 0.05956416464891041
$ws.Range("R3").Value = 0.623574144486692
$ws.Range("U3").Value = 6.81
$ws.Range("V3").Value = 0.0824455205811138
$ws.Range("W3").Value = 0.1230889235569423
$ws.Range("X3").Value = 0.06302613697497982
$ws.Range("Y3").Value = 0.06006278658196247
$ws.Range("Z3").Value = 2.112960161371659
$ws.Range("AA3").Value = 0.1354434839716967
$ws.Range("AB3").Value = 0.06302613697497982
$ws.Range("AC3").Value = 0.07241734699671691
$ws.Range("AG3").Value = -6.81
$ws.Range("AJ3").Value = -0.08985354268373137
$ws.Range("AK3").Value = -0.1150532184490623
$ws.Range("AL3").Value = 0.304
$ws.Range("AM3").Value = 0.304
$ws.Range("AO3").Value = 33.22368421052632
$ws.Range("AP3").Value = -0.5675
$ws.Range("AQ3").Value = 33.22368421052632

# Row 4
$ws.Range("D4").Value = 0.13
$ws.Range("E4").Value = 0.245
$ws.Range("G4").Value = 0.1306755260243632
$ws.Range("H4").Value = 0.1306755260243632
$ws.Range("I4").Value = 0.1604079598706648
$ws.Range("J4").Value = 0.1222568775230472
$ws.Range("K4").Value = 11.1
$ws.Range("L4").Value = 0.1229235880398671
$ws.Range("M4").Value = 3.7
$ws.Range("N4").Value = 0.03162393162393162
$ws.Range("O4").Value = 0.3333333333333334
$ws.Range("P4").Value = 3.7
$ws.Range("Q4").Value = 0.03162393162393162
$ws.Range("R4").Value = 0.3333333333333334
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 4.43
$ws.Range("V4").Value = 0.03786324786324786
$ws.Range("W4").Value = 0.1183368869936034
$ws.Range("X4").Value = 0.06345807899599111
$ws.Range("Y4").Value = 0.0548788079976123
$ws.Range("Z4").Value = 0.989960006304103
$ws.Range("AA4").Value = 0.1210294192434357
$ws.Range("AB4").Value = 0.06331852767618502
$ws.Range("AC4").Value = 0.05771089156725072
$ws.Range("AE4").Value = 1.175806118394853
$ws.Range("AF4").Value = 1.175806118394853
$ws.Range("AG4").Value = -3.254193881605147
$ws.Range("AH4").Value = 0.009949634845027984
$ws.Range("AI4").Value = 0.01146280163801702
$ws.Range("AJ4").Value = -0.02860935266675192
$ws.Range("AK4").Value = -0.03315672885379902
$ws.Range("AL4").Value = 0.003
$ws.Range("AM4").Value = 0.003
$ws.Range("AO4").Value = 4866.666666666666
$ws.Range("AP4").Value = -0.2166573822639911
$ws.Range("AQ4").Value = 4866.666666666666

# Row 5
$ws.Range("B5").Value = "Post & Telecommunication Joint Stock Insurance Corporation (HNX:PTI)"
$ws.Range("D5").Value = 0.267
$ws.Range("E5").Value = 0.235
$ws.Range("G5").Value = 0.0251860119047619
$ws.Range("H5").Value = 0.0251860119047619
$ws.Range("I5").Value = 0.04575892857142858
$ws.Range("J5").Value = 0.03570696721311475
$ws.Range("K5").Value = 9.59
$ws.Range("L5").Value = 0.03567708333333333
$ws.Range("M5").Value = 3.4572
$ws.Range("N5").Value = 0.04542969776609724
$ws.Range("O5").Value = 0.3605005213764338
$ws.Range("P5").Value = 3.4572
$ws.Range("Q5").Value = 0.04542969776609724
$ws.Range("R5").Value = 0.3605005213764338
$ws.Range("U5").Value = 6.2
$ws.Range("V5").Value = 0.08147174770039423
$ws.Range("W5").Value = 0.123741935483871
$ws.Range("X5").Value = 0.06884352683928775
$ws.Range("Y5").Value = 0.05489840864458322
$ws.Range("Z5").Value = 3.294036910860028
$ws.Range("AA5").Value = 0.1176200679748688
$ws.Range("AB5").Value = 0.06613080953924051
$ws.Range("AC5").Value = 0.05148925843562832
$ws.Range("AD5").Value = 10.3
$ws.Range("AF5").Value = 10.3
$ws.Range("AG5").Value = 4.100000000000001
$ws.Range("AH5").Value = 0.119212962962963
$ws.Range("AI5").Value = 0.1098081023454158
$ws.Range("AJ5").Value = 0.05112219451371573
$ws.Range("AK5").Value = 0.04680365296803654
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = -0.193
$ws.Range("AN5").Value = 0.7518248175182483
$ws.Range("AP5").Value = 0.2992700729927008
$ws.Range("AQ5").Value = -63.73056994818653

# Remove AO5 cell entirely (column dropped for this row)
$ws.Range("AO5").ClearContents()

Write-Output "Done"